$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77, shifting existing rows 77:111 down to 78:112
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new weekly record
$ws.Range("A77").Value = 10
$ws.Range("B77").Value = "Vega Modelo de Temuco"
$ws.Range("C77").Value = "La Araucanía"
$ws.Range("D77").Value = 44784
$ws.Range("E77").Value = 9
$ws.Range("F77").Value = 100112035
$ws.Range("G77").Value = "Bruselas (repollito)"
$ws.Range("H77").Value = "Sin especificar"
$ws.Range("I77").Value = "Primera"
$ws.Range("J77").Value = 180
$ws.Range("K77").Value = 24000
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = 24444
$ws.Range("N77").Value = "$/malla 10 kilos"
$ws.Range("O77").Value = "Región Metropolitana"
$ws.Range("P77").Value = 2444
$ws.Range("Q77").Value = 10
$ws.Range("R77").Value = "Hortaliza"
